# Insert a new weekly price record for "Vega Monumental Concepción -
# Arándano (blue)" above the current row 68, shifting the existing
# rows 68-85 down to 69-86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68; Excel shifts rows 68:85 down to 69:86
# and copies the row-68 formatting (incl. the date style on column D) into
# the freshly inserted row, matching the original workbook's pattern.
$ws.Rows.Item(68).Insert()

# Populate the new row 68 with the new record.
$ws.Range("A68").Value = 11
$ws.Range("B68").Value = "Vega Monumental Concepción"
$ws.Range("C68").Value = "Bíobío"
$ws.Range("D68").Value = 44588
$ws.Range("E68").Value = 8
$ws.Range("F68").Value = "Fruta"
$ws.Range("G68").Value = 100101
$ws.Range("H68").Value = "Berries"
$ws.Range("I68").Value = 100101001
$ws.Range("J68").Value = "Arándano (blue)"
$ws.Range("K68").Value = "Sin especificar"
$ws.Range("L68").Value = "Primera"
$ws.Range("M68").Value = 150
$ws.Range("N68").Value = 3000
$ws.Range("O68").Value = 3500
$ws.Range("P68").Value = 3267
$ws.Range("Q68").Value = "$/bandeja 2 kilos"
$ws.Range("R68").Value = "Provincia de Linares"
$ws.Range("S68").Value = 1634
$ws.Range("T68").Value = 2
